# Added new triangulation2 script
# Update the computed "Actual"/"Estimated" values for the three mic pairs
# (columns B:G, rows 2-11) with the results produced by the new script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.9112048919340322, 1.133786848072562, 0.1126721327937182, 0.9070294784580499, 1.139435185830596, 0.1133786848072562),
    @(-0.4402577480568004, -0.36281179138322, 0.1028323251606915, -0.4308390022675737, -0.3662529298129745, 0.1133786848072562),
    @(-2.031350183375844, -2.448979591836735, -1.253243669098878, -2.040816326530612, -2.444198566859693, -1.26984126984127),
    @(-0.002692718193406324, -0.2040816326530612, -0.1855199464135275, 0, -0.1878768028333454, -0.18140589569161),
    @(1.541215953830037, 1.927437641723356, 0.1231218379758774, 1.541950113378685, 1.945556573671575, 0.1133786848072562),
    @(-0.3802581138671877, 0.18140589569161, 0.7815463054709451, -0.3854875283446712, 0.1838007771989899, 0.7709750566893424),
    @(0.5145133491967829, 0.453514739229025, -0.04511345882168556, 0.5215419501133787, 0.449486283507921, -0.04535147392290249),
    @(0.4907005132262501, 1.224489795918367, 0.4794669711284663, 0.4988662131519275, 1.220539079278846, 0.4761904761904762),
    @(1.165374245474661, 2.244897959183673, 0.3733966858293393, 1.179138321995465, 2.245223706833127, 0.3854875283446712),
    @(0.8777519492556363, 0.4308390022675737, -0.2530379391380345, 0.8843537414965986, 0.4345807984913791, -0.2494331065759637)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
